# Rename the "Test Cases" sheet to "Test_Cases" and make it the active
# sheet (instead of "Data"), so the workbook opens on it next time.
$wb = $excel.ActiveWorkbook

$testCasesSheet = $wb.Worksheets.Item("Test Cases")
$testCasesSheet.Name = "Test_Cases"

# Activating this sheet moves Excel's active/selected tab here, clearing
# the tabSelected flag that currently sits on "Data".
$testCasesSheet.Activate()
